$d = $word.ActiveDocument

# 1. Delete the team-info table (Name/ID + 4 team members)
$d.Tables(1).Delete()

# 2. Delete the "Team Information:" and "Team ID: 2" paragraphs, together
#    with the now-orphaned empty paragraph that used to trail the table.
$pStart = $d.Paragraphs(57)
$pEnd = $d.Paragraphs(59)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rng.Delete()

# 3. Strip the pStyle (ListParagraph) and ind (left=996) overrides from the
#    remaining empty paragraph, while preserving its spacing/jc/rPr.
$p = $d.Paragraphs(56)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)
